$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 119
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value2 = [Math]::Round([double]$current, 0)
    }
}
